# CSV file import implemented
# - Renumber the CueDis_* labels in column A (rows 2-5): 6/7/8/9 -> 19/20/21/22
# - Remove the trailing "CueDis_10" row (row 6) entirely - the CSV import now
#   produces one fewer row, so A6 goes back to being empty.
# - Column widths were recalculated slightly narrower/wider after the import.
# - Selection moved from B10 to A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the existing CueDis_* identifiers (descriptions in column B are untouched)
$ws.Range("A2").Value = "CueDis_19"
$ws.Range("A3").Value = "CueDis_20"
$ws.Range("A4").Value = "CueDis_21"
$ws.Range("A5").Value = "CueDis_22"

# Row 6 (previously "CueDis_10") no longer exists after the import - clear it fully
$ws.Range("A6").Clear()

# Updated column widths (A, B, and the default width for the remaining columns)
$ws.Columns.Item(1).ColumnWidth = 34.166666666666664
$ws.Columns.Item(2).ColumnWidth = 38.33333333333333
$ws.Columns.Item(3).ColumnWidth = 7.666666666666667

# Active cell/selection moved to A10
[void]$ws.Range("A10").Select()
